$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22)
$newSerial = 45221

$ws.Range("C2").Value = $newSerial
$ws.Range("C3").Value = $newSerial
$ws.Range("C4").Value = $newSerial
$ws.Range("C5").Value = $newSerial
$ws.Range("C6").Value = $newSerial
